# Apply the "polish the remaining accel filter logic shift" edit:
# - Fill in rows 18-20 on the vbRcSignals sheet with aeb/fcw/dbs request signal mappings
# - Switch the active sheet/selection back to vbRcSignals (from params)

$wb = $excel.ActiveWorkbook

$wsSignals = $wb.Worksheets.Item("vbRcSignals")
$wsParams  = $wb.Worksheets.Item("params")

# --- Fill in the new signal rows on vbRcSignals ---
# New unique strings must appear in this exact order so the shared-string
# table is built the same way as the target workbook:
#   aebRequest, fcw_request, fcwRequest, aeb_request, dbs_request, dbsRequest
$wsSignals.Range("A18").Value = "aebRequest"
$wsSignals.Range("B18").Value = "fcw_request"
$wsSignals.Range("C18").Value = "fcwRequest"

$wsSignals.Range("A19").Value = "aebRequest"
$wsSignals.Range("B19").Value = "aeb_request"
$wsSignals.Range("C19").Value = "aebRequest"

$wsSignals.Range("A20").Value = "aebRequest"
$wsSignals.Range("B20").Value = "dbs_request"
$wsSignals.Range("C20").Value = "dbsRequest"

# --- Update the view state: vbRcSignals becomes the active tab again ---
# Update the (inactive) params sheet's remembered selection first ...
$wsParams.Range("G22").Select() | Out-Null

# ... then switch back to vbRcSignals and leave it selected/active, matching
# the new cursor position (frozen header row stays frozen; the cursor moves
# from B23 to C21).
$wsSignals.Activate()
$wsSignals.Range("C21").Select() | Out-Null

$wb.Save()
